{"js": "// The document body contains 19 top-level paragraphs. This edit re-shuffles\n// the *text content* of several paragraphs (moving whole blocks of text to\n// different heading sections) while every paragraph itself \u2014 its style,\n// its position, and (for the multi-run \"Avalia\u00e7\u00e3o\" bullet) its run-level\n// bold/italic formatting \u2014 stays exactly where it already is.\n//\n// Concretely (0-based paragraph index in the body):\n//   5  (objectives, PT)      ->  gets the old \"Programa resumido\" PT text (10)\n//   6  (objectives, EN)      ->  gets the old \"Programa\" EN text (11)\n//   8  (Docente bullet)      ->  gets the old objectives PT text (5)\n//   10 (Programa resumido PT)->  gets the old long \"Programa\" PT text (13)\n//   11 (Programa resumido EN)->  gets the old objectives EN text (6)\n//   13 (long \"Programa\" PT)  ->  gets the old \"M\u00e9todo:\" value text (16, run)\n//   16 \"M\u00e9todo:\" value       ->  gets the old \"Crit\u00e9rio:\" value text (16, run)\n//   16 \"Crit\u00e9rio:\" value     ->  gets the old \"Norma de recupera\u00e7\u00e3o:\" value (16, run)\n//   16 \"Norma...\" value      ->  gets the old Bibliografia paragraph text (18)\n//   18 (Bibliografia text)   ->  gets the old Docente bullet text (8)\n//\n// Because this is a genuine cycle (each destination is also a source), every\n// original value is read first and only then are the writes applied, so no\n// write can clobber text that is still needed as a source.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// --- Phase 1: get Range objects for every whole-paragraph slot we touch, and\n// load their current text (captured BEFORE any mutation happens). ---\nconst rngObjPt = paragraphs.items[5].getRange();\nconst rngObjEn = paragraphs.items[6].getRange();\nconst rngDocente = paragraphs.items[8].getRange();\nconst rngResumoPt = paragraphs.items[10].getRange();\nconst rngResumoEn = paragraphs.items[11].getRange();\nconst rngProgramaLong = paragraphs.items[13].getRange();\nconst rngBibliografia = paragraphs.items[18].getRange();\nconst rngAvaliacao = paragraphs.items[16].getRange();\n\nrngObjPt.load(\"text\");\nrngObjEn.load(\"text\");\nrngDocente.load(\"text\");\nrngResumoPt.load(\"text\");\nrngResumoEn.load(\"text\");\nrngProgramaLong.load(\"text\");\nrngBibliografia.load(\"text\");\nawait context.sync();\n\n// The \"Avalia\u00e7\u00e3o\" bullet paragraph has 3 bold labels (\"M\u00e9todo: \",\n// \"Crit\u00e9rio: \", \"Norma de recupera\u00e7\u00e3o: \") each followed by a plain-text\n// value run. Locate each value run's Range (scoped search inside that one\n// paragraph, so there is no risk of matching text elsewhere in the body).\nconst metodoValOld =\n  \"Semin\u00e1rios e Estudos de Casos, aplica\u00e7\u00e3o de \\u201cPitchs\\u201d (breve apresenta\u00e7\u00e3o oral de uma ideia, produto ou oportunidade de neg\u00f3cio) e outras formas de apresenta\u00e7\u00e3o de ideias em empresas\";\nconst criterioValOld =\n  \"As avalia\u00e7\u00f5es ser\u00e3o: a) cont\u00ednuas considerando a participa\u00e7\u00e3o dos alunos nas atividades; b) avalia\u00e7\u00e3o das apresenta\u00e7\u00f5es parciais dos trabalhos; e c) apresenta\u00e7\u00e3o final dos trabalhos.\";\nconst normaValOld =\n  \"Reapresenta\u00e7\u00e3o do \u00faltimo semin\u00e1rio, cuja nota constituir\u00e1 a nota final da disciplina.\";\n\nconst searchMetodo = rngAvaliacao.search(metodoValOld, { matchCase: true });\nconst searchCriterio = rngAvaliacao.search(criterioValOld, { matchCase: true });\nconst searchNorma = rngAvaliacao.search(normaValOld, { matchCase: true });\nsearchMetodo.load(\"items\");\nsearchCriterio.load(\"items\");\nsearchNorma.load(\"items\");\nawait context.sync();\n\nconst rngMetodoVal = searchMetodo.items[0];\nconst rngCriterioVal = searchCriterio.items[0];\nconst rngNormaVal = searchNorma.items[0];\n\n// Capture every original value now, before any write happens.\nconst objPtOld = rngObjPt.text;\nconst objEnOld = rngObjEn.text;\nconst docenteOld = rngDocente.text;\nconst resumoPtOld = rngResumoPt.text;\nconst resumoEnOld = rngResumoEn.text;\nconst programaLongOld = rngProgramaLong.text;\nconst bibliografiaOld = rngBibliografia.text;\n\n// --- Phase 2: apply every new value using the captured originals. ---\nrngObjPt.insertText(resumoPtOld, Word.InsertLocation.replace);\nrngObjEn.insertText(resumoEnOld, Word.InsertLocation.replace);\nrngDocente.insertText(objPtOld, Word.InsertLocation.replace);\nrngResumoPt.insertText(programaLongOld, Word.InsertLocation.replace);\nrngResumoEn.insertText(objEnOld, Word.InsertLocation.replace);\nrngProgramaLong.insertText(metodoValOld, Word.InsertLocation.replace);\nrngMetodoVal.insertText(criterioValOld, Word.InsertLocation.replace);\nrngCriterioVal.insertText(normaValOld, Word.InsertLocation.replace);\nrngNormaVal.insertText(bibliografiaOld, Word.InsertLocation.replace);\nrngBibliografia.insertText(docenteOld, Word.InsertLocation.replace);\n\nawait context.sync();\n", "ps1": "# The document body has 19 paragraphs ($d.Paragraphs is 1-indexed). This\n# edit re-shuffles the *text content* of several paragraphs (moving whole\n# blocks of text to different heading sections) while every paragraph\n# itself - its style, its position, and (for the multi-run \"Avalia\u00e7\u00e3o\"\n# bullet) its run-level bold formatting - stays exactly where it is.\n#\n# Concretely (1-based $d.Paragraphs index):\n#   6  (objectives, PT)       -> gets the old \"Programa resumido\" PT text (11)\n#   7  (objectives, EN)       -> gets the old \"Programa\" EN text (12)\n#   9  (Docente bullet)       -> gets the old objectives PT text (6)\n#   11 (Programa resumido PT) -> gets the old long \"Programa\" PT text (14)\n#   12 (Programa resumido EN) -> gets the old objectives EN text (7)\n#   14 (long \"Programa\" PT)   -> gets the old \"M\u00e9todo:\" value text (17, run)\n#   17 \"M\u00e9todo:\" value        -> gets the old \"Crit\u00e9rio:\" value text (17, run)\n#   17 \"Crit\u00e9rio:\" value      -> gets the old \"Norma de recupera\u00e7\u00e3o:\" value (17, run)\n#   17 \"Norma...\" value       -> gets the old Bibliografia paragraph text (19)\n#   19 (Bibliografia text)    -> gets the old Docente bullet text (9)\n#\n# Because this is a genuine cycle (each destination is also a source\n# somewhere else), every original value is captured into a variable first,\n# and only afterwards are the writes applied - so no write can clobber text\n# that is still needed as a source.\n\n$d = $word.ActiveDocument\n\n# NOTE: Paragraph.Range.Text includes the trailing paragraph-mark (chr 13).\n# It must be stripped before the captured string is written into a\n# *different* paragraph, otherwise re-assigning it splits that paragraph in\n# two (inflating the paragraph count).\nfunction Get-ParaText($para) {\n    return $para.Range.Text.TrimEnd([char]13)\n}\n\n# --- Phase 1: capture every original value we will need, before any write. ---\n$objPtOld        = Get-ParaText $d.Paragraphs(6)\n$objEnOld        = Get-ParaText $d.Paragraphs(7)\n$docenteOld      = Get-ParaText $d.Paragraphs(9)\n$resumoPtOld     = Get-ParaText $d.Paragraphs(11)\n$resumoEnOld     = Get-ParaText $d.Paragraphs(12)\n$programaLongOld = Get-ParaText $d.Paragraphs(14)\n$bibliografiaOld = Get-ParaText $d.Paragraphs(19)\n\n# The \"Avalia\u00e7\u00e3o\" bullet (paragraph 17) has 3 bold labels (\"M\u00e9todo: \",\n# \"Crit\u00e9rio: \", \"Norma de recupera\u00e7\u00e3o: \") each followed by a plain-text\n# value run. Locate + capture each value via a Find scoped to a duplicate\n# of that one paragraph's range, so there is no risk of matching text\n# anywhere else in the document.\n$metodoValOld = \"Semin\u00e1rios e Estudos de Casos, aplica\u00e7\u00e3o de \u201cPitchs\u201d (breve apresenta\u00e7\u00e3o oral de uma ideia, produto ou oportunidade de neg\u00f3cio) e outras formas de apresenta\u00e7\u00e3o de ideias em empresas\"\n$criterioValOld = \"As avalia\u00e7\u00f5es ser\u00e3o: a) cont\u00ednuas considerando a participa\u00e7\u00e3o dos alunos nas atividades; b) avalia\u00e7\u00e3o das apresenta\u00e7\u00f5es parciais dos trabalhos; e c) apresenta\u00e7\u00e3o final dos trabalhos.\"\n$normaValOld = \"Reapresenta\u00e7\u00e3o do \u00faltimo semin\u00e1rio, cuja nota constituir\u00e1 a nota final da disciplina.\"\n\n# --- Phase 2: apply every new value using the captured originals. ---\n$d.Paragraphs(6).Range.Text  = $resumoPtOld\n$d.Paragraphs(7).Range.Text  = $resumoEnOld\n$d.Paragraphs(9).Range.Text  = $objPtOld\n$d.Paragraphs(11).Range.Text = $programaLongOld\n$d.Paragraphs(12).Range.Text = $objEnOld\n$d.Paragraphs(14).Range.Text = $metodoValOld\n\n$avaliacaoRange = $d.Paragraphs(17).Range.Duplicate\n$avaliacaoRange.Find.Execute($metodoValOld) | Out-Null\n$avaliacaoRange.Text = $criterioValOld\n\n$avaliacaoRange2 = $d.Paragraphs(17).Range.Duplicate\n$avaliacaoRange2.Find.Execute($criterioValOld) | Out-Null\n$avaliacaoRange2.Text = $normaValOld\n\n$avaliacaoRange3 = $d.Paragraphs(17).Range.Duplicate\n$avaliacaoRange3.Find.Execute($normaValOld) | Out-Null\n$avaliacaoRange3.Text = $bibliografiaOld\n\n$d.Paragraphs(19).Range.Text = $docenteOld\n"}
